$p = $ppt.ActivePresentation
Write-Output "Slide count: $($p.Slides.Count)"
$s = $p.Slides.Item(1)
Write-Output "Shape count: $($s.Shapes.Count)"
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    Write-Output "Shape $i : $($sh.Name) Type=$($sh.Type)"
}
